$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3, shifting existing rows 3..130 down to 4..131
$ws.Rows.Item(3).Insert()

# Copy static formatting/values from row 2 into new row 3 for columns that are
# constant across the whole sheet (A,B,C,E,F,G,H,I,N,O,Q,R), then set the
# varying values (D,J,K,L,M,P) for this new record.
$ws.Range("A2:R2").Copy()
$ws.Range("A3:R3").PasteSpecial(-4104)

$ws.Cells.Item(3, 4).Value = 44860
$ws.Cells.Item(3, 10).Value = 1000
$ws.Cells.Item(3, 11).Value = 4000
$ws.Cells.Item(3, 12).Value = 4500
$ws.Cells.Item(3, 13).Value = 4250
$ws.Cells.Item(3, 16).Value = 4250
